$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2023" data column (column P), mirroring column O's formatting ---

# Row 3 is a blank border-only row; just carry over O3's border style.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# Rows 4-14: year header + the new data points, using each row's column-O
# cell as the format donor so the new column inherits the right numFmt/
# font/border (and thus reuses the existing style indexes rather than
# manufacturing new ones).
$newColumnValues = @{
  4  = 2023
  5  = 48.2
  6  = 8.6767564891727478
  7  = 12.226605469730881
  8  = 78.520866131691164
  9  = 59.466452648968115
  10 = 26.635270208942913
  11 = 8.166450559693871
  12 = 74.601894583630667
  13 = 99.168063426054971
  14 = 70.956108992253434
}

foreach ($row in 4..14) {
  $ws.Range("O$row").Copy()
  $ws.Range("P$row").PasteSpecial(-4122)
  $ws.Range("P$row").Value = $newColumnValues[$row]
}

# --- Row 14's D/E cells: previously blank, now show a right-aligned "-" ---
$ws.Range("D14").Value = "-"
$ws.Range("D14").HorizontalAlignment = -4152
$ws.Range("E14").Value = "-"
$ws.Range("E14").HorizontalAlignment = -4152

# --- Row height tweaks (data rows grow from default/13.5pt to 15pt; the
#     thin spacer row grows from 7.5pt to 13.5pt) ---
foreach ($row in 4..14) {
  $ws.Rows.Item($row).RowHeight = 15
}
$ws.Rows.Item(15).RowHeight = 13.5

# --- Clear the lingering P8 selection left over in the sheet view ---
$ws.Range("A1").Select() | Out-Null
